# Update "想去人数" (column F) values across the 展览 / 演出 / 全部类型 sheets
# to match the newly generated gh-pages output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    2  = 5938
    4  = 1158
    5  = 1085
    8  = 55
    9  = 627
    10 = 70
    11 = 39
    13 = 2127
    14 = 1538
    15 = 1189
    16 = 305
    18 = 468
    19 = 692
    23 = 525
    24 = 3932
    26 = 137
    28 = 178
    29 = 62
    30 = 570
    35 = 340
    36 = 885
    37 = 116
    38 = 80
    39 = 97
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# --- Sheet "演出" (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$sheet2Updates = @{
    3 = 765
    8 = 1
    9 = 14
}
foreach ($row in $sheet2Updates.Keys) {
    $ws2.Cells.Item($row, 6).Value = $sheet2Updates[$row]
}

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    2  = 5938
    4  = 1158
    6  = 765
    7  = 1085
    12 = 55
    13 = 627
    14 = 70
    15 = 39
    18 = 2127
    19 = 1538
    20 = 1189
    21 = 305
    23 = 468
    25 = 692
    29 = 525
    30 = 3932
    32 = 137
    34 = 178
    35 = 62
    36 = 570
    41 = 340
    42 = 885
    43 = 116
    44 = 80
    45 = 97
    47 = 1
    48 = 14
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
